# Auto update Excel log
# Appends newly logged sensor readings to the PIR, Humidity and mmWave
# sheets of the SeniorConnect master log workbook.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [object[]]$Rows
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $endRow = $StartRow + $Rows.Count - 1

    # Force column A (Date) to remain plain text so strings like
    # "2026-01-30" are not reinterpreted as date serial numbers.
    $ws.Range("A$StartRow`:A$endRow").NumberFormat = "@"
    # Force column E (Value) to remain plain text so strings like
    # "87.8%" are not reinterpreted as numeric percentages.
    $ws.Range("E$StartRow`:E$endRow").NumberFormat = "@"

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $rowValues = $Rows[$i]
        $r = $StartRow + $i
        for ($j = 0; $j -lt $rowValues.Count; $j++) {
            $ws.Cells.Item($r, $j + 1).Value = $rowValues[$j]
        }
    }
}

# PIR sheet: new motion-sensor log rows 103-117
$pirRows = @(
        @("2026-01-30", "15:45:59", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:00", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:04", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:09", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:14", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:19", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:20", "15:00", "Living Room", "RECOVERY_DETECTION", "Inactive"),
        @("2026-01-30", "15:46:24", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:29", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:34", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:39", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:44", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:49", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:54", "15:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-01-30", "15:46:59", "15:00", "Bathroom", "No Motion", "Inactive")
)
Add-LogRows "PIR" 103 $pirRows

# Humidity sheet: new humidity-sensor log rows 59-68
$humidityRows = @(
        @("2026-01-30", "15:46:00", "15:00", "Bathroom", "87.8%", "Active"),
        @("2026-01-30", "15:46:04", "15:00", "Bathroom", "87.8%", "Active"),
        @("2026-01-30", "15:46:09", "15:00", "Bathroom", "86.9%", "Active"),
        @("2026-01-30", "15:46:14", "15:00", "Bathroom", "87.8%", "Active"),
        @("2026-01-30", "15:46:24", "15:00", "Bathroom", "87.8%", "Active"),
        @("2026-01-30", "15:46:29", "15:00", "Bathroom", "86.9%", "Active"),
        @("2026-01-30", "15:46:34", "15:00", "Bathroom", "87.7%", "Active"),
        @("2026-01-30", "15:46:44", "15:00", "Bathroom", "87.7%", "Active"),
        @("2026-01-30", "15:46:49", "15:00", "Bathroom", "86.9%", "Active"),
        @("2026-01-30", "15:46:55", "15:00", "Bathroom", "87.7%", "Active")
)
Add-LogRows "Humidity" 59 $humidityRows

# mmWave sheet: new presence-sensor log rows 10-13
$mmwaveRows = @(
        @("2026-01-30", "15:46:20", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
        @("2026-01-30", "15:46:30", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
        @("2026-01-30", "15:46:40", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
        @("2026-01-30", "15:46:51", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
)
Add-LogRows "mmWave" 10 $mmwaveRows

Write-Host "Log rows appended."
